{"js": "// Replace each old two-digit multiplication equation with its new value.\n// Every \"from\" string is unique in the document, so a direct search +\n// insertText(\"Replace\") on the first (only) hit is unambiguous.\nconst replacements = [\n  { from: \"23\u00d793=2139\", to: \"56\u00d749=2744\" },\n  { from: \"54\u00d764=3456\", to: \"53\u00d737=1961\" },\n  { from: \"59\u00d763=3717\", to: \"13\u00d742=546\" },\n  { from: \"57\u00d768=3876\", to: \"71\u00d770=4970\" },\n  { from: \"92\u00d715=1380\", to: \"96\u00d719=1824\" },\n  { from: \"69\u00d753=3657\", to: \"60\u00d777=4620\" },\n  { from: \"46\u00d791=4186\", to: \"19\u00d782=1558\" },\n  { from: \"32\u00d735=1120\", to: \"46\u00d794=4324\" },\n  { from: \"45\u00d733=1485\", to: \"42\u00d788=3696\" },\n  { from: \"78\u00d767=5226\", to: \"49\u00d782=4018\" },\n  { from: \"69\u00d754=3726\", to: \"84\u00d717=1428\" },\n  { from: \"28\u00d722=616\", to: \"59\u00d755=3245\" },\n  { from: \"77\u00d740=3080\", to: \"38\u00d746=1748\" },\n  { from: \"28\u00d741=1148\", to: \"86\u00d772=6192\" },\n  { from: \"31\u00d778=2418\", to: \"16\u00d748=768\" },\n  { from: \"57\u00d785=4845\", to: \"11\u00d789=979\" },\n  { from: \"19\u00d768=1292\", to: \"17\u00d768=1156\" },\n  { from: \"66\u00d759=3894\", to: \"41\u00d745=1845\" },\n  { from: \"79\u00d759=4661\", to: \"75\u00d764=4800\" },\n  { from: \"65\u00d739=2535\", to: \"48\u00d754=2592\" },\n  { from: \"34\u00d767=2278\", to: \"93\u00d741=3813\" },\n  { from: \"58\u00d717=986\", to: \"61\u00d722=1342\" },\n  { from: \"92\u00d752=4784\", to: \"70\u00d765=4550\" },\n  { from: \"33\u00d783=2739\", to: \"13\u00d788=1144\" },\n  { from: \"81\u00d770=5670\", to: \"84\u00d789=7476\" }\n];\n\nconst body = context.document.body;\n\nfor (const { from, to } of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${from}`);\n  }\n\n  results.items[0].insertText(to, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace each old two-digit multiplication equation with its new value.\n# Every \"Find\" string is unique in the document, so Find/Replace (wdReplaceAll,\n# but only ever matching a single occurrence) is unambiguous for each pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Find = \"23\u00d793=2139\"; Replace = \"56\u00d749=2744\" }\n  @{ Find = \"54\u00d764=3456\"; Replace = \"53\u00d737=1961\" }\n  @{ Find = \"59\u00d763=3717\"; Replace = \"13\u00d742=546\" }\n  @{ Find = \"57\u00d768=3876\"; Replace = \"71\u00d770=4970\" }\n  @{ Find = \"92\u00d715=1380\"; Replace = \"96\u00d719=1824\" }\n  @{ Find = \"69\u00d753=3657\"; Replace = \"60\u00d777=4620\" }\n  @{ Find = \"46\u00d791=4186\"; Replace = \"19\u00d782=1558\" }\n  @{ Find = \"32\u00d735=1120\"; Replace = \"46\u00d794=4324\" }\n  @{ Find = \"45\u00d733=1485\"; Replace = \"42\u00d788=3696\" }\n  @{ Find = \"78\u00d767=5226\"; Replace = \"49\u00d782=4018\" }\n  @{ Find = \"69\u00d754=3726\"; Replace = \"84\u00d717=1428\" }\n  @{ Find = \"28\u00d722=616\"; Replace = \"59\u00d755=3245\" }\n  @{ Find = \"77\u00d740=3080\"; Replace = \"38\u00d746=1748\" }\n  @{ Find = \"28\u00d741=1148\"; Replace = \"86\u00d772=6192\" }\n  @{ Find = \"31\u00d778=2418\"; Replace = \"16\u00d748=768\" }\n  @{ Find = \"57\u00d785=4845\"; Replace = \"11\u00d789=979\" }\n  @{ Find = \"19\u00d768=1292\"; Replace = \"17\u00d768=1156\" }\n  @{ Find = \"66\u00d759=3894\"; Replace = \"41\u00d745=1845\" }\n  @{ Find = \"79\u00d759=4661\"; Replace = \"75\u00d764=4800\" }\n  @{ Find = \"65\u00d739=2535\"; Replace = \"48\u00d754=2592\" }\n  @{ Find = \"34\u00d767=2278\"; Replace = \"93\u00d741=3813\" }\n  @{ Find = \"58\u00d717=986\"; Replace = \"61\u00d722=1342\" }\n  @{ Find = \"92\u00d752=4784\"; Replace = \"70\u00d765=4550\" }\n  @{ Find = \"33\u00d783=2739\"; Replace = \"13\u00d788=1144\" }\n  @{ Find = \"81\u00d770=5670\"; Replace = \"84\u00d789=7476\" }\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute(\n    $pair.Find,    # FindText\n    $true,         # MatchCase\n    $false,        # MatchWholeWord\n    $false,        # MatchWildcards\n    $false,        # MatchSoundsLike\n    $false,        # MatchAllWordForms\n    $true,         # Forward\n    1,             # Wrap            (wdFindContinue)\n    $false,        # Format\n    $pair.Replace, # ReplaceWith\n    2              # Replace         (wdReplaceAll)\n  )\n}\n"}
